$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.85
$ws.Range("T2").Value = 1.36

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("V3").Value = 1.18

# Row 4
$ws.Range("G4").Value = 1.6
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2.25
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.33
$ws.Range("V4").Value = 1.25
$ws.Range("AB4").Value = 7
$ws.Range("AI4").Value = 19

# Row 5
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("V5").Value = 1.14
